$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 onto the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 6

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
